$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.8809986114502
$ws.Range("B3").Value = 7.95662784576416
$ws.Range("B4").Value = 5.955420017242432
$ws.Range("B5").Value = 6.767120838165283
$ws.Range("B6").Value = 7.603854179382324
$ws.Range("B7").Value = 11.47830581665039
$ws.Range("B8").Value = 13.53664112091064
$ws.Range("B9").Value = 10.49045848846436
$ws.Range("B10").Value = 4.737107753753662
$ws.Range("B11").Value = 6.079762935638428
$ws.Range("B12").Value = 4.378615856170654
$ws.Range("B13").Value = 7.389666080474854
$ws.Range("B14").Value = 10.40564441680908
$ws.Range("B15").Value = 10.42655563354492
$ws.Range("B16").Value = 10.09898853302002
$ws.Range("B17").Value = 6.089284420013428
$ws.Range("B18").Value = 4.424148559570312
$ws.Range("B19").Value = 13.29117012023926
$ws.Range("B20").Value = 19.33496475219727
$ws.Range("B21").Value = 6.577813625335693
$ws.Range("B22").Value = 5.151782989501953
$ws.Range("B23").Value = 4.395680904388428
$ws.Range("B24").Value = 6.069608688354492
$ws.Range("B25").Value = 6.620934963226318
$ws.Range("B26").Value = 7.318948745727539
$ws.Range("B27").Value = 11.81636333465576
$ws.Range("B28").Value = 8.639409065246582
$ws.Range("B29").Value = 11.24605464935303
$ws.Range("B30").Value = 14.32633876800537
$ws.Range("B31").Value = 11.50711727142334
$ws.Range("B32").Value = 15.95943927764893
$ws.Range("B33").Value = 6.54820442199707
$ws.Range("B34").Value = 19.80498123168945
$ws.Range("B35").Value = 25.43582534790039
$ws.Range("B36").Value = 13.52932548522949
$ws.Range("B37").Value = 7.902271747589111
$ws.Range("B38").Value = 6.021476745605469
$ws.Range("B39").Value = 6.838600158691406
$ws.Range("B40").Value = 7.617969036102295
$ws.Range("B41").Value = 11.15285015106201
$ws.Range("B42").Value = 13.29177665710449
$ws.Range("B43").Value = 10.23600673675537
$ws.Range("B44").Value = 4.950057506561279
$ws.Range("B45").Value = 6.18470287322998
$ws.Range("B46").Value = 4.31477689743042
$ws.Range("B47").Value = 7.165870189666748
$ws.Range("B48").Value = 10.13792991638184
$ws.Range("B49").Value = 9.974912643432617
$ws.Range("B50").Value = 9.893525123596191
$ws.Range("B51").Value = 5.997543811798096
$ws.Range("B52").Value = 4.566448211669922
$ws.Range("B53").Value = 12.97424793243408
$ws.Range("B54").Value = 19.20149993896484
$ws.Range("B55").Value = 6.523789882659912
$ws.Range("B56").Value = 5.140477657318115
$ws.Range("B57").Value = 4.459598064422607
$ws.Range("B58").Value = 6.025604724884033
$ws.Range("B59").Value = 6.619915008544922
$ws.Range("B60").Value = 7.2478346824646
$ws.Range("B61").Value = 11.04336071014404
$ws.Range("B62").Value = 8.525527000427246
$ws.Range("B63").Value = 11.04471969604492
$ws.Range("B64").Value = 13.73667907714844
$ws.Range("B65").Value = 11.18331527709961
$ws.Range("B66").Value = 15.89662742614746
$ws.Range("B67").Value = 6.532844543457031
$ws.Range("B68").Value = 19.44545364379883
$ws.Range("B69").Value = 24.73703193664551
$ws.Range("B70").Value = 13.4672269821167
$ws.Range("B71").Value = 7.831560611724854
$ws.Range("B72").Value = 6.095869064331055
$ws.Range("B73").Value = 6.864828586578369
$ws.Range("B74").Value = 7.593636989593506
$ws.Range("B75").Value = 11.05364608764648
$ws.Range("B76").Value = 13.05564117431641
$ws.Range("B77").Value = 10.1266508102417
$ws.Range("B78").Value = 5.177108287811279
$ws.Range("B79").Value = 6.270837306976318
$ws.Range("B80").Value = 4.252202033996582
$ws.Range("B81").Value = 6.946003913879395
$ws.Range("B82").Value = 9.864738464355469
$ws.Range("B83").Value = 9.410311698913574
$ws.Range("B84").Value = 9.71759033203125
$ws.Range("B85").Value = 5.906689167022705
$ws.Range("B86").Value = 4.805705070495605
$ws.Range("B87").Value = 12.64804077148438
$ws.Range("B88").Value = 19.15567207336426
$ws.Range("B89").Value = 6.473563671112061
$ws.Range("B90").Value = 5.091256618499756
$ws.Range("B91").Value = 4.499393939971924
$ws.Range("B92").Value = 6.052968978881836
$ws.Range("B93").Value = 6.607266902923584
$ws.Range("B94").Value = 7.175524711608887
$ws.Range("B95").Value = 10.60240840911865
$ws.Range("B96").Value = 8.40052318572998
$ws.Range("B97").Value = 10.70535087585449
$ws.Range("B98").Value = 13.24186706542969
$ws.Range("B99").Value = 11.0146427154541
$ws.Range("B100").Value = 15.82266807556152
$ws.Range("B101").Value = 6.48128604888916
$ws.Range("B102").Value = 19.08477592468262
$ws.Range("B103").Value = 23.96520614624023
$ws.Range("B104").Value = 13.74034023284912
$ws.Range("B105").Value = 7.841311931610107
$ws.Range("B106").Value = 6.168068885803223
$ws.Range("B107").Value = 6.835846900939941
$ws.Range("B108").Value = 7.528514862060547
$ws.Range("B109").Value = 10.78582572937012
$ws.Range("B110").Value = 12.74608516693115
$ws.Range("B111").Value = 10.11339473724365
$ws.Range("B112").Value = 5.427739143371582
$ws.Range("B113").Value = 6.331255912780762
$ws.Range("B114").Value = 4.204568862915039
$ws.Range("B115").Value = 6.737100601196289
$ws.Range("B116").Value = 9.57958984375
$ws.Range("B117").Value = 8.886554718017578
$ws.Range("B118").Value = 9.530409812927246
$ws.Range("B119").Value = 5.814504623413086
$ws.Range("B120").Value = 4.861351013183594
$ws.Range("B121").Value = 12.26574897766113
$ws.Range("B122").Value = 18.90688514709473
$ws.Range("B123").Value = 6.423656940460205
$ws.Range("B124").Value = 5.126580238342285
$ws.Range("B125").Value = 4.519901275634766
$ws.Range("B126").Value = 6.120871543884277
$ws.Range("B127").Value = 6.603997707366943
$ws.Range("B128").Value = 7.106751918792725
$ws.Range("B129").Value = 10.43985939025879
$ws.Range("B130").Value = 8.278443336486816
$ws.Range("B131").Value = 10.38094711303711
$ws.Range("B132").Value = 12.79089641571045
$ws.Range("B133").Value = 10.76794719696045
$ws.Range("B134").Value = 15.74698162078857
$ws.Range("B135").Value = 6.446595668792725
$ws.Range("B136").Value = 18.74972724914551
$ws.Range("B137").Value = 23.27029037475586
